# Scheduled market-data refresh: update cached price/profit figures
# across the per-class Leve worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 116.111115
$ws.Range("I11").Value = 116.111115
$ws.Range("K11").Value = 116.111115
$ws.Range("M11").Value = 23.888885
$ws.Range("H55").Value = 826.93335
$ws.Range("J55").Value = 934.4545000000001
$ws.Range("L55").Value = 934.4545000000001
$ws.Range("N55").Value = -1362.4545
$ws.Range("H63").Value = 75577.75
$ws.Range("H66").Value = 75577.75
$ws.Range("H93").Value = 69696
$ws.Range("J93").Value = 69696
$ws.Range("L93").Value = 69696
$ws.Range("N93").Value = -74688
$ws.Range("H97").Value = 3428.8572
$ws.Range("J97").Value = 3428.8572
$ws.Range("L97").Value = 10286.5716
$ws.Range("N97").Value = -11278.5716
$ws.Range("H99").Value = 957.44446
$ws.Range("I99").Value = 960
$ws.Range("J99").Value = 948.5
$ws.Range("K99").Value = 2880
$ws.Range("L99").Value = 2845.5
$ws.Range("M99").Value = -1382
$ws.Range("N99").Value = -5841.5
$ws.Range("H101").Value = 1127.1111
$ws.Range("I101").Value = 540.6667
$ws.Range("K101").Value = 1622.0001
$ws.Range("M101").Value = -0.00009999999997489795
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H110").Value = 56567.332
$ws.Range("J110").Value = 56567.332
$ws.Range("L110").Value = 56567.332
$ws.Range("N110").Value = -64747.332
$ws.Range("H112").Value = 2501367.2
$ws.Range("J112").Value = 3334493
$ws.Range("L112").Value = 10003479
$ws.Range("N112").Value = -10005695
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
$ws.Range("H115").Value = 967.75
$ws.Range("I115").Value = 298.66666
$ws.Range("K115").Value = 895.9999799999999
$ws.Range("M115").Value = 671.0000200000001
$ws.Range("H117").Value = 95531.664
$ws.Range("J117").Value = 95531.664
$ws.Range("L117").Value = 95531.664
$ws.Range("N117").Value = -104709.664
$ws.Range("H118").Value = 1274.4
$ws.Range("J118").Value = 2066.6667
$ws.Range("L118").Value = 6200.000100000001
$ws.Range("N118").Value = -9514.000100000001
$ws.Range("H129").Value = 1299.8684
$ws.Range("I129").Value = 528.2857
$ws.Range("K129").Value = 1584.8571
$ws.Range("M129").Value = 3415.1429
$ws.Range("H135").Value = 5753
$ws.Range("I135").Value = 2664.4443
$ws.Range("J135").Value = 9227.625
$ws.Range("K135").Value = 23979.9987
$ws.Range("L135").Value = 83048.625
$ws.Range("M135").Value = -21444.9987
$ws.Range("N135").Value = -88118.625
$ws.Range("H138").Value = 226487.9
$ws.Range("J138").Value = 282625.34
$ws.Range("L138").Value = 847876.02
$ws.Range("N138").Value = -858156.02

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 36734.25
$ws.Range("J46").Value = 37707.2
$ws.Range("L46").Value = 37707.2
$ws.Range("N46").Value = -38345.2
$ws.Range("H74").Value = 4710.207
$ws.Range("I74").Value = 2849.75
$ws.Range("J74").Value = 7000
$ws.Range("K74").Value = 2849.75
$ws.Range("L74").Value = 7000
$ws.Range("M74").Value = -1975.75
$ws.Range("N74").Value = -8748
$ws.Range("H77").Value = 4710.207
$ws.Range("I77").Value = 2849.75
$ws.Range("J77").Value = 7000
$ws.Range("K77").Value = 14248.75
$ws.Range("L77").Value = 35000
$ws.Range("M77").Value = -9880.75
$ws.Range("N77").Value = -43736
$ws.Range("H97").Value = 6088.625
$ws.Range("I97").Value = 1331.1666
$ws.Range("J97").Value = 20361
$ws.Range("K97").Value = 1331.1666
$ws.Range("L97").Value = 20361
$ws.Range("M97").Value = -835.1666
$ws.Range("N97").Value = -21353
$ws.Range("H132").Value = 3358.0981
$ws.Range("I132").Value = 2557.8809
$ws.Range("J132").Value = 7092.4443
$ws.Range("K132").Value = 7673.6427
$ws.Range("L132").Value = 21277.3329
$ws.Range("M132").Value = -5143.6427
$ws.Range("N132").Value = -26337.3329
$ws.Range("H139").Value = 86989.39999999999
$ws.Range("J139").Value = 86989.39999999999
$ws.Range("L139").Value = 86989.39999999999
$ws.Range("N139").Value = -97269.39999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 34485.5
$ws.Range("J26").Value = 42000
$ws.Range("L26").Value = 42000
$ws.Range("N26").Value = -42584
$ws.Range("H54").Value = 6109.25
$ws.Range("I54").Value = 5630.3335
$ws.Range("K54").Value = 5630.3335
$ws.Range("M54").Value = -5146.3335
$ws.Range("H96").Value = 8448.333000000001
$ws.Range("I96").Value = 8448.333000000001
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 8448.333000000001
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -5702.333000000001
$ws.Range("N96").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2463.5293
$ws.Range("I31").Value = 1043.2273
$ws.Range("K31").Value = 1043.2273
$ws.Range("M31").Value = -748.2273
$ws.Range("H34").Value = 2463.5293
$ws.Range("I34").Value = 1043.2273
$ws.Range("K34").Value = 1043.2273
$ws.Range("M34").Value = -841.2273
$ws.Range("H76").Value = 4951
$ws.Range("I76").Value = 4951
$ws.Range("K76").Value = 4951
$ws.Range("M76").Value = -4636
$ws.Range("H79").Value = 4951
$ws.Range("I79").Value = 4951
$ws.Range("K79").Value = 4951
$ws.Range("M79").Value = -3859
$ws.Range("H99").Value = 7091.923
$ws.Range("I99").Value = 6971.5557
$ws.Range("K99").Value = 6971.5557
$ws.Range("M99").Value = -5473.5557
$ws.Range("H105").Value = 4508.636
$ws.Range("I105").Value = 4759.5
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 4759.5
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -3012.5
$ws.Range("N105").Value = -5494
$ws.Range("H122").Value = 4424.778
$ws.Range("I122").Value = 2637.1667
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 7911.500100000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5461.500100000001
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 7091.923
$ws.Range("I126").Value = 6971.5557
$ws.Range("K126").Value = 20914.6671
$ws.Range("M126").Value = -18444.6671

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 9962
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 9962
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 29886
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -30340

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 23966.137
$ws.Range("I102").Value = 967.82855
$ws.Range("J102").Value = 113404
$ws.Range("K102").Value = 967.82855
$ws.Range("L102").Value = 113404
$ws.Range("M102").Value = 654.17145
$ws.Range("N102").Value = -116648
$ws.Range("H132").Value = 3566.585
$ws.Range("J132").Value = 1977.2
$ws.Range("L132").Value = 5931.6
$ws.Range("N132").Value = -10991.6

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1863.8889
$ws.Range("I16").Value = 1708.4445
$ws.Range("K16").Value = 1708.4445
$ws.Range("M16").Value = -1538.4445
$ws.Range("H56").Value = 30042.727
$ws.Range("I56").Value = 26992.25
$ws.Range("J56").Value = 38177.332
$ws.Range("K56").Value = 26992.25
$ws.Range("L56").Value = 38177.332
$ws.Range("M56").Value = -26301.25
$ws.Range("N56").Value = -39559.332
$ws.Range("H68").Value = 2549.375
$ws.Range("I68").Value = 2549.375
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2549.375
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1800.375
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2549.375
$ws.Range("I71").Value = 2549.375
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 12746.875
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9002.875
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 1120.8182
$ws.Range("I93").Value = 1120.8182
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1120.8182
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 127.1818000000001
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 3585.6428
$ws.Range("I100").Value = 3219.9
$ws.Range("K100").Value = 3219.9
$ws.Range("M100").Value = -2678.9

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10010
$ws.Range("I20").Value = 10010
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10010
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9770
$ws.Range("N20").ClearContents()
$ws.Range("H51").Value = 14561.5
$ws.Range("I51").Value = 7116
$ws.Range("K51").Value = 7116
$ws.Range("M51").Value = -6606
$ws.Range("H96").Value = 1564.375
$ws.Range("I96").Value = 1473.5714
$ws.Range("J96").Value = 2200
$ws.Range("K96").Value = 1473.5714
$ws.Range("L96").Value = 2200
$ws.Range("M96").Value = -100.5714
$ws.Range("N96").Value = -4946
$ws.Range("H100").Value = 1087.3
$ws.Range("I100").Value = 1128.6666
$ws.Range("K100").Value = 2257.3332
$ws.Range("M100").Value = -1716.3332
$ws.Range("H119").Value = 65000
$ws.Range("J119").Value = 65000
$ws.Range("L119").Value = 65000
$ws.Range("N119").Value = -74676
$ws.Range("H132").Value = 1448.8
$ws.Range("I132").Value = 1351.7646
$ws.Range("K132").Value = 4055.2938
$ws.Range("M132").Value = -1525.2938
$ws.Range("H140").Value = 99997
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99997
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99997
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -110357

